# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Column A holds a numeric rank.
#
# Some Price cells parse as a clean number (e.g. "213.60"); a bare
# Range.Value assignment would auto-convert those to numeric cells, same as
# typing them straight into Excel. To keep them as plain text (matching the
# source data, which stores these as strings like "68.57" / "1.00"), we
# momentarily force a Text number format, assign the string, then
# ClearFormats() so the cell's style/format reverts to the sheet default
# (General) without losing the now-text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = '@'
    $r.Value = $text
    $r.ClearFormats()
}

$ws.Range('D2').Value = '27.915.10'
$ws.Range('D3').Value = '1.644.22'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.60'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.56'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0617'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('E11').Value = '  -1.35%  '
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').Value = '1.649.39'
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.574'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.44%  '
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.89'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').Value = '27.896.18'
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.36'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').Value = '0.0₃0725'
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.64'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.89'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +4.48%  '
$ws.Range('E24').Value = '  +2.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.50'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.93'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('E27').Value = '  +0.89%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D33').Value = '1.425.26'
$ws.Range('E33').Value = '  -2.93%  '
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('E35').Value = '  +1.58%  '
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.886'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.925'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.77%  '
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.04'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.18%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '68.57'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('E45').Value = '  +2.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.81'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.76%  '
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = '1.785.82'
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '88.93'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.72'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.77%  '
